$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column G ("PIB marítimo") entirely; columns H:L shift left to G:K.
$ws.Columns("G:G").Delete()

# Update the "transporte, correos y almacenamiento" footnote to fold in the
# maritime-transport wording that used to live in its own note.
$ws.Range("B35").Value = "  El concepto PIB transporte, correos y almacenamiento, incluye transporte maritimo, transporte por ductos, transporte turístico, servicios relacionados con el transporte, y servicios de almacenamiento."
